$d = $word.ActiveDocument

# Locate the two list-paragraphs that discuss "other possible tables and/or
# graphs": the first one asks to tabulate/chart campaign counts per
# Category/Sub-Category, the second asks to tabulate campaign duration.
# The edit removes the first paragraph entirely and keeps the second one
# (together with its bold paragraph-mark formatting), moving the existing
# "_GoBack" bookmark from the middle of that paragraph to its very start.

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Create a table and chart that will analyze the*Sub-Category*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    # Remove the whole paragraph (its text and its paragraph mark); the
    # following paragraph (the duration/chart one) slides up and keeps its
    # own formatting (including the bold paragraph mark).
    $target.Range.Delete()
}

# Find the paragraph that now starts with "Create a table to calculate a
# campaign duration..." so we can relocate the bookmark to its beginning.
$dest = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Create a table to c*pledged is higher*") {
        $dest = $p
        break
    }
}

if ($dest -ne $null -and $d.Bookmarks.Exists("_GoBack")) {
    $startPos = $dest.Range.Start
    $bm = $d.Bookmarks.Item("_GoBack")
    $bm.Delete()
    $r = $d.Range($startPos, $startPos)
    $d.Bookmarks.Add("_GoBack", $r)
}
